# FrontendTimeline.xlsx - add a new timeline entry (Day 7) and move the
# "Total hours Spent" summary block down one row to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the "Total hours Spent" row (row 11),
# which pushes that summary block (and the following merged cells)
# down by one row, while leaving row 10 free for the new entry.
$ws.Rows(11).Insert()

# Fill in the new timeline entry in row 10.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "19/5/2024"
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = "Added categories in FE with design + Search Bar"

# Match the centered alignment used by the rest of the data rows.
$ws.Cells.Item(10, 1).HorizontalAlignment = -4108
$ws.Cells.Item(10, 3).HorizontalAlignment = -4108
$ws.Cells.Item(10, 4).HorizontalAlignment = -4108

# Extend the totals formula (now on row 12 after the insert) to include
# the newly added row.
$ws.Cells.Item(12, 4).Formula = "=SUM(C4:C10)"

# The formula edit nudges the row's auto height; re-fit it back to the
# sheet's default so no custom row height sticks around.
$ws.Rows(12).AutoFit()

# Match the saved selection location from the source file.
$ws.Range("D15").Select()
